$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.328.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.565.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.560.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.38%  "
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.670"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -11.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.141.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.566.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.296.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.76%  "
$ws.Range("E21").Value = "  -6.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "616.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.113"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "62.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.389"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0757"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.131"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.038.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.24%  "
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("E46").Value = "  -7.85%  "
$ws.Range("E47").Value = "  -7.57%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.11%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
